$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 326, shifting rows 326-346 down to 327-347.
$ws.Rows.Item(326).Insert()

# Populate the newly inserted row 326 with the new data record.
$ws.Range("A326").Value = 10
$ws.Range("B326").Value = "Vega Modelo de Temuco"
$ws.Range("C326").Value = "La Araucanía"
$ws.Range("D326").Value = 45021
$ws.Range("E326").Value = 9
$ws.Range("F326").Value = "Fruta"
$ws.Range("G326").Value = 100103
$ws.Range("H326").Value = "Frutos de hueso (carozo)"
$ws.Range("I326").Value = 100103002
$ws.Range("J326").Value = "Ciruela"
$ws.Range("K326").Value = "Angeleno"
$ws.Range("L326").Value = "Primera"
$ws.Range("M326").Value = 100
$ws.Range("N326").Value = 14000
$ws.Range("O326").Value = 14000
$ws.Range("P326").Value = 14000
$ws.Range("Q326").Value = "$/bandeja 18 kilos granel"
$ws.Range("R326").Value = "Región de O'Higgins"
$ws.Range("S326").Value = 778
$ws.Range("T326").Value = 18
